# Generate Report for Handoff
# Update the "Latest Handoff Datetime" (column D, row 4) for the
# 0302ff9b-fb7a-4eb0-999d-c98d43afa1eb file on both the zh-cn and de-de
# language sheets to reflect the newly generated handoff report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-03-11 08:07:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-03-11 08:08:08"
